$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: force a genuine run split at a sub-range by toggling a character
# property off/on (no net formatting change) around the mutation. Direct
# Range.Text assignment alone gets silently re-merged into the neighbouring
# run when formatting is identical, so we nudge the bold flag to keep the
# boundary.
# -----------------------------------------------------------------------
function Split-SetText($rng, [string]$text) {
    $rng.Bold = $true
    $rng.Text = $text
    $rng.Bold = $false
}

function Split-Touch($rng) {
    $rng.Bold = $true
    $rng.Bold = $false
}

# Insert `text` right after the (collapsed) range `rng` and force the newly
# inserted text to live in its own run (rather than being silently absorbed
# into the preceding run because the formatting matches).
function Insert-Split($rng, [string]$text) {
    $pos = $rng.End
    $rng.InsertAfter($text)
    $rNew = $d.Range($pos, $pos + $text.Length)
    $rNew.Bold = $true
    $rNew.Bold = $false
    return $rNew
}

# 1) Motivation paragraph: "python and matlab" -> "Python and MATLAB"
#    with the exact run split pattern: "P" | "ython and " | "MATLAB"
$r = $d.Content
$r.Find.Execute("python and matlab so", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $r.Start

$rP = $d.Range($s, $s + 1)
Split-SetText $rP "P"

$rYthonAnd = $d.Range($s + 1, $s + 1 + 10)   # "ython and "
Split-Touch $rYthonAnd

$rMatlab = $d.Range($s + 11, $s + 11 + 6)    # "matlab"
Split-SetText $rMatlab "MATLAB"

# 2) ENU table cell: "East north up (m,m,m,lat,long,alt)." ->
#    "East north up (m,m,m). sensorloc holds the origin."
#    Run split pattern: "E" | "ast north up (m,m,m" | ")." | " sensorloc holds the origin."
$r = $d.Content
$r.Find.Execute("East north up (m,m,m,lat,long,alt).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $r.Start
$e = $r.End

$rE = $d.Range($s, $s + 1)
Split-Touch $rE

$rAst = $d.Range($s + 1, $s + 20)    # "ast north up (m,m,m"
Split-Touch $rAst

$rParen = $d.Range($s + 20, $e)  # ",lat,long,alt)."
Split-SetText $rParen ")."

$rSensor = $d.Range($rParen.End, $rParen.End)
Insert-Split $rSensor " sensorloc holds the origin." | Out-Null

# 3) Cartesian row: append " Pretty much the same as ENU but in km."
$r = $d.Content
$r.Find.Execute("Local Cartesian grid (km,km,km).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rEnd = $d.Range($r.End, $r.End)
Insert-Split $rEnd " Pretty much the same as ENU but in km." | Out-Null

# 4) "python setup.m" -> "python setup.py" (first occurrence, standalone command)
$r = $d.Content
$r.Find.Execute("python setup.m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $r.Start
$rM = $d.Range($s + 13, $s + 14)   # the "m"
Split-SetText $rM "py"

# Move the _GoBack bookmark to the end of this paragraph's content
$rBookmark = $d.Range($rM.End, $rM.End)
$d.Bookmarks.Add("_GoBack", $rBookmark)

# 5) "python setup.m develop" -> "python setup.py develop" (second occurrence)
$r = $d.Content
$r.Find.Execute("python setup.m develop", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $r.Start
$rM2 = $d.Range($s + 13, $s + 14)  # the "m"
Split-SetText $rM2 "py"
$rDevelop = $d.Range($rM2.End, $rM2.End + 8)  # " develop"
Split-Touch $rDevelop
